$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph 1 ("Image Replacement in global paragraphs", Titre1):
# split the single run into three runs: "Image " / "Replacement" /
# " in global paragraphs" (the word "Replacement" becomes its own run).
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1Start = $p1.Range.Start
$rMid = $d.Range($p1Start + 6, $p1Start + 17)
# Toggling a character property forces Word to split the run at these
# boundaries; toggling it back off keeps the visible formatting intact.
$rMid.Bold = 1
$rMid.Bold = 0

# ------------------------------------------------------------------
# Paragraph 2 ("This paragraph is untouched."): drop the explicit
# "Textkoerper" paragraph style so it falls back to the default style.
# ------------------------------------------------------------------
$d.Paragraphs(2).Style = "Normal"

# ------------------------------------------------------------------
# Paragraph 3 (Mona Lisa inserted in the sentence): drop the explicit
# paragraph style and collapse the "${" / "monalisa" / "}." runs
# (and their surrounding proofErr spell-check markers) back into a
# single run.
# ------------------------------------------------------------------
$d.Paragraphs(3).Style = "Normal"
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$r3.End = $r3.End - 1
$r3.Text = "PLACEHOLDER_TEXT_3"
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$r3.End = $r3.End - 1
$r3.Text = 'In this paragraph, an image of Mona Lisa is inserted: ${monalisa}.'

# ------------------------------------------------------------------
# Paragraph 4 (image in the middle of the sentence): same treatment.
# ------------------------------------------------------------------
$d.Paragraphs(4).Style = "Normal"
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$r4.End = $r4.End - 1
$r4.Text = "PLACEHOLDER_TEXT_4"
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$r4.End = $r4.End - 1
$r4.Text = 'This paragraph has the image ${monalisa} in the middle.'
